# Edit script for Temp24monthsStatisticsCalculator.xlsx
# Replaces the McKenzie/Hayden statistics (2010-18) with NSantiam/Detroit
# statistics (2019-20): new observed/simulated site labels, new H/I data
# series, refreshed year labels, and a D-column index that is now computed
# via an incrementing formula instead of hard-coded numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics calculator")

# --- Site header labels (H3: site name, I3: observation file path) ---
$ws.Range("H3").Value = " USGS_14178000_temp_NO SANTIAM R BLW BOULDER CRK  NR DETROIT  OR_23780591"
$ws.Range("I3").Value = " Obs:..\Observations\NSantiam\USGS_14178000_temp_NO SANTIAM R BLW BOULDER CRK  NR DETROIT  OR_23780591.csv"

# --- Update H/I temperature data (simulated & observed) for NSantiam site ---
$ws.Cells.Item(4, 8).Value = 5.263011
$ws.Cells.Item(4, 9).Value = 4.110979
$ws.Cells.Item(5, 8).Value = 3.968922
$ws.Cells.Item(5, 9).Value = 2.857552
$ws.Cells.Item(6, 8).Value = 5.397714
$ws.Cells.Item(6, 9).Value = 3.784271
$ws.Cells.Item(7, 8).Value = 6.790036
$ws.Cells.Item(7, 9).Value = 5.474094
$ws.Cells.Item(8, 8).Value = 8.040923
$ws.Cells.Item(8, 9).Value = 8.54385
$ws.Cells.Item(9, 8).Value = 9.633939
$ws.Cells.Item(9, 9).Value = 11.242885
$ws.Cells.Item(10, 8).Value = 10.283904
$ws.Cells.Item(10, 9).Value = 12.688961
$ws.Cells.Item(11, 8).Value = 10.65511
$ws.Cells.Item(11, 9).Value = 13.284344
$ws.Cells.Item(12, 8).Value = 8.720452
$ws.Cells.Item(12, 9).Value = 11.13861
$ws.Cells.Item(13, 8).Value = 6.043349
$ws.Cells.Item(13, 9).Value = 6.833771
$ws.Cells.Item(14, 8).Value = 6.040865
$ws.Cells.Item(14, 9).Value = 5.036613
$ws.Cells.Item(15, 8).Value = 4.365733
$ws.Cells.Item(15, 9).Value = 4.099933
$ws.Cells.Item(16, 8).Value = 4.691846
$ws.Cells.Item(16, 9).Value = 3.99773
$ws.Cells.Item(17, 8).Value = 4.965358
$ws.Cells.Item(17, 9).Value = 4.055855
$ws.Cells.Item(18, 8).Value = 5.225219
$ws.Cells.Item(18, 9).Value = 4.478296
$ws.Cells.Item(19, 8).Value = 6.8973
$ws.Cells.Item(19, 9).Value = 5.938178
$ws.Cells.Item(20, 8).Value = 7.927648
$ws.Cells.Item(20, 9).Value = 8.001192
$ws.Cells.Item(21, 8).Value = 9.074625
$ws.Cells.Item(21, 9).Value = 10.241529
$ws.Cells.Item(22, 8).Value = 10.537362
$ws.Cells.Item(22, 9).Value = 12.795609
$ws.Cells.Item(23, 8).Value = 10.470796
$ws.Cells.Item(23, 9).Value = 12.995021
$ws.Cells.Item(24, 8).Value = 9.504089
$ws.Cells.Item(24, 9).Value = 10.865839
$ws.Cells.Item(25, 8).Value = 7.537137
$ws.Cells.Item(25, 9).Value = 8.217125
$ws.Cells.Item(26, 8).Value = 5.07178
$ws.Cells.Item(26, 9).Value = 5.324552
$ws.Cells.Item(27, 8).Value = 4.829082
$ws.Cells.Item(27, 9).Value = 4.357759

# --- Update year column E: 2010->2019 (rows 4-15), 2018->2020 (rows 16-27) ---
$ws.Cells.Item(4, 5).Value = 2019
$ws.Cells.Item(5, 5).Value = 2019
$ws.Cells.Item(6, 5).Value = 2019
$ws.Cells.Item(7, 5).Value = 2019
$ws.Cells.Item(8, 5).Value = 2019
$ws.Cells.Item(9, 5).Value = 2019
$ws.Cells.Item(10, 5).Value = 2019
$ws.Cells.Item(11, 5).Value = 2019
$ws.Cells.Item(12, 5).Value = 2019
$ws.Cells.Item(13, 5).Value = 2019
$ws.Cells.Item(14, 5).Value = 2019
$ws.Cells.Item(15, 5).Value = 2019
$ws.Cells.Item(16, 5).Value = 2020
$ws.Cells.Item(17, 5).Value = 2020
$ws.Cells.Item(18, 5).Value = 2020
$ws.Cells.Item(19, 5).Value = 2020
$ws.Cells.Item(20, 5).Value = 2020
$ws.Cells.Item(21, 5).Value = 2020
$ws.Cells.Item(22, 5).Value = 2020
$ws.Cells.Item(23, 5).Value = 2020
$ws.Cells.Item(24, 5).Value = 2020
$ws.Cells.Item(25, 5).Value = 2020
$ws.Cells.Item(26, 5).Value = 2020
$ws.Cells.Item(27, 5).Value = 2020

# --- Convert D16:D27 from static values to incrementing formulas ---
$ws.Cells.Item(16, 4).Formula = "=D15+1"
$ws.Cells.Item(17, 4).Formula = "=D16+1"
$ws.Cells.Item(18, 4).Formula = "=D17+1"
$ws.Cells.Item(19, 4).Formula = "=D18+1"
$ws.Cells.Item(20, 4).Formula = "=D19+1"
$ws.Cells.Item(21, 4).Formula = "=D20+1"
$ws.Cells.Item(22, 4).Formula = "=D21+1"
$ws.Cells.Item(23, 4).Formula = "=D22+1"
$ws.Cells.Item(24, 4).Formula = "=D23+1"
$ws.Cells.Item(25, 4).Formula = "=D24+1"
$ws.Cells.Item(26, 4).Formula = "=D25+1"
$ws.Cells.Item(27, 4).Formula = "=D26+1"

# --- NSE (B6) now displayed with 4 decimal places instead of 3 ---
$ws.Range("B6").NumberFormat = "0.0000"

# --- Selection moves from H3:I27 to H4:I27 ---
$ws.Activate()
$null = $ws.Range("H4:I27").Select()
